# Auto-generated edit script applying value updates to Titan_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 134 (ALC)
$ws.Range("H134").Value = 45778.75
$ws.Range("J134").Value = 45778.75
$ws.Range("L134").Value = 45778.75
$ws.Range("N134").Value = -55918.75

# Row 136 (ALC)
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 132 (ARM)
$ws.Range("H132").Value = 2745.125
$ws.Range("I132").Value = 2262.1155
$ws.Range("K132").Value = 6786.3465
$ws.Range("M132").Value = -4256.3465

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (BSM)
$ws.Range("H105").Value = 2897.8667
$ws.Range("I105").Value = 2970.5417
$ws.Range("J105").Value = 2814.8096
$ws.Range("K105").Value = 2970.5417
$ws.Range("L105").Value = 2814.8096
$ws.Range("M105").Value = -1223.5417
$ws.Range("N105").Value = -6308.809600000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 1417.2
$ws.Range("I31").Value = 1089.2858
$ws.Range("J31").Value = 2182.3333
$ws.Range("K31").Value = 1089.2858
$ws.Range("L31").Value = 2182.3333
$ws.Range("M31").Value = -794.2858000000001
$ws.Range("N31").Value = -2772.3333

# Row 34 (CRP)
$ws.Range("H34").Value = 1417.2
$ws.Range("I34").Value = 1089.2858
$ws.Range("J34").Value = 2182.3333
$ws.Range("K34").Value = 1089.2858
$ws.Range("L34").Value = 2182.3333
$ws.Range("M34").Value = -887.2858000000001
$ws.Range("N34").Value = -2586.3333

# Row 86 (CRP)
$ws.Range("H86").Value = 41668860
$ws.Range("I86").Value = 83335090
$ws.Range("J86").Value = 2633
$ws.Range("K86").Value = 83335090
$ws.Range("L86").Value = 2633
$ws.Range("M86").Value = -83333967
$ws.Range("N86").Value = -4879

# Row 89 (CRP)
$ws.Range("H89").Value = 41668860
$ws.Range("I89").Value = 83335090
$ws.Range("J89").Value = 2633
$ws.Range("K89").Value = 416675450
$ws.Range("L89").Value = 13165
$ws.Range("M89").Value = -416669834
$ws.Range("N89").Value = -24397

# Row 99 (CRP)
$ws.Range("H99").Value = 7814012.5
$ws.Range("J99").Value = 1900
$ws.Range("L99").Value = 1900
$ws.Range("N99").Value = -4896

# Row 122 (CRP)
$ws.Range("H122").Value = 1083.3334
$ws.Range("I122").Value = 1080
$ws.Range("J122").Value = 1100
$ws.Range("K122").Value = 3240
$ws.Range("L122").Value = 3300
$ws.Range("M122").Value = -790
$ws.Range("N122").Value = -8200

# Row 126 (CRP)
$ws.Range("H126").Value = 7814012.5
$ws.Range("J126").Value = 1900
$ws.Range("L126").Value = 5700
$ws.Range("N126").Value = -10640

$ws = $wb.Worksheets.Item("CUL")
# Row 108 (CUL)
$ws.Range("H108").Value = 3287.5
$ws.Range("J108").Value = 6000
$ws.Range("L108").Value = 18000
$ws.Range("N108").Value = -23760

# Row 137 (CUL)
$ws.Range("H137").Value = 4813056
$ws.Range("I137").Value = 10002114
$ws.Range("J137").Value = 95730.27
$ws.Range("K137").Value = 30006342
$ws.Range("L137").Value = 287190.81
$ws.Range("M137").Value = -30001242
$ws.Range("N137").Value = -297390.81

# Row 138 (CUL)
$ws.Range("H138").Value = 952
$ws.Range("I138").Value = 952
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2856
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = 2284

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 3325
$ws.Range("I80").Value = 3250
$ws.Range("J80").Value = 3400
$ws.Range("K80").Value = 3250
$ws.Range("L80").Value = 3400
$ws.Range("M80").Value = -2252
$ws.Range("N80").Value = -5396

# Row 83 (GSM)
$ws.Range("H83").Value = 3325
$ws.Range("I83").Value = 3250
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 16250
$ws.Range("L83").Value = 17000
$ws.Range("M83").Value = -11258
$ws.Range("N83").Value = -26984

# Row 102 (GSM)
$ws.Range("H102").Value = 1731.125
$ws.Range("I102").Value = 1287.5
$ws.Range("J102").Value = 2174.75
$ws.Range("K102").Value = 1287.5
$ws.Range("L102").Value = 2174.75
$ws.Range("M102").Value = 334.5
$ws.Range("N102").Value = -5418.75

# Row 122 (GSM)
$ws.Range("H122").Value = 5557055.5
$ws.Range("I122").Value = 11111111
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 33333333
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -33330883
$ws.Range("N122").Value = -13900

# Row 132 (GSM)
$ws.Range("H132").Value = 2713.52
$ws.Range("I132").Value = 1957.9286
$ws.Range("J132").Value = 3675.182
$ws.Range("K132").Value = 5873.7858
$ws.Range("L132").Value = 11025.546
$ws.Range("M132").Value = -3343.7858
$ws.Range("N132").Value = -16085.546

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 3500
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

# Row 40 (LTW)
$ws.Range("H40").Value = 4095.652
$ws.Range("I40").Value = 2666.6667
$ws.Range("J40").Value = 4310
$ws.Range("K40").Value = 2666.6667
$ws.Range("L40").Value = 4310
$ws.Range("N40").Value = -4582
$ws.Range("M40").Value = -2530.6667

# Row 82 (LTW)
$ws.Range("H82").Value = 27971.053
$ws.Range("I82").Value = 37160.715
$ws.Range("J82").Value = 2240
$ws.Range("K82").Value = 37160.715
$ws.Range("L82").Value = 2240
$ws.Range("M82").Value = -36799.715
$ws.Range("N82").Value = -2962

# Row 85 (LTW)
$ws.Range("H85").Value = 27971.053
$ws.Range("I85").Value = 37160.715
$ws.Range("J85").Value = 2240
$ws.Range("K85").Value = 37160.715
$ws.Range("L85").Value = 2240
$ws.Range("M85").Value = -35912.715
$ws.Range("N85").Value = -4736

# Row 126 (LTW)
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (WVR)
$ws.Range("H113").Value = 1733.3334
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 5700
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = -3530
$ws.Range("N113").Value = -8540

# Row 126 (WVR)
$ws.Range("H126").Value = 63937.75
$ws.Range("I126").Value = 112267.11
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 336801.33
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -334331.33
$ws.Range("N126").Value = -10340

# Row 136 (WVR)
$ws.Range("H136").Value = 17598012
$ws.Range("I136").Value = 22289990
$ws.Range("J136").Value = 3099.5
$ws.Range("K136").Value = 66869970
$ws.Range("L136").Value = 9298.5
$ws.Range("M136").Value = -66867420
$ws.Range("N136").Value = -14398.5
